# Scheduled runner update: refresh market-price / profit columns (H-N)
# across the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 333337.66
$ws.Cells.Item(11, 9).Value = 333337.66
$ws.Cells.Item(11, 11).Value = 333337.66
$ws.Cells.Item(11, 13).Value = -333197.66

$ws.Cells.Item(17, 8).Value = 676.65717
$ws.Cells.Item(17, 10).Value = 687.82355
$ws.Cells.Item(17, 12).Value = 2063.47065
$ws.Cells.Item(17, 14).Value = -2399.47065

$ws.Cells.Item(51, 8).Value = 20483.166
$ws.Cells.Item(51, 10).Value = 5399.6665
$ws.Cells.Item(51, 12).Value = 5399.6665
$ws.Cells.Item(51, 14).Value = -6367.6665

$ws.Cells.Item(55, 8).Value = 638
$ws.Cells.Item(55, 9).Value = 643.8
$ws.Cells.Item(55, 10).Value = 630.75
$ws.Cells.Item(55, 11).Value = 643.8
$ws.Cells.Item(55, 12).Value = 630.75
$ws.Cells.Item(55, 13).Value = -429.8
$ws.Cells.Item(55, 14).Value = -1058.75

$ws.Cells.Item(76, 8).Value = 4960
$ws.Cells.Item(76, 9).Value = 4800
$ws.Cells.Item(76, 10).Value = 5000
$ws.Cells.Item(76, 11).Value = 4800
$ws.Cells.Item(76, 12).Value = 5000
$ws.Cells.Item(76, 13).Value = -4485
$ws.Cells.Item(76, 14).Value = -5630

$ws.Cells.Item(79, 8).Value = 4960
$ws.Cells.Item(79, 9).Value = 4800
$ws.Cells.Item(79, 10).Value = 5000
$ws.Cells.Item(79, 11).Value = 4800
$ws.Cells.Item(79, 12).Value = 5000
$ws.Cells.Item(79, 13).Value = -3708
$ws.Cells.Item(79, 14).Value = -7184

$ws.Cells.Item(127, 8).Value = 1851.7949
$ws.Cells.Item(127, 10).Value = 2103.9092
$ws.Cells.Item(127, 12).Value = 6311.7276
$ws.Cells.Item(127, 14).Value = -16231.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 575.2
$ws.Cells.Item(22, 9).Value = 575.2
$ws.Cells.Item(22, 11).Value = 575.2
$ws.Cells.Item(22, 13).Value = -276.2

$ws.Cells.Item(32, 8).Value = 25189.125
$ws.Cells.Item(32, 9).Value = 4162.7256
$ws.Cells.Item(32, 10).Value = 155552.8
$ws.Cells.Item(32, 11).Value = 4162.7256
$ws.Cells.Item(32, 12).Value = 155552.8
$ws.Cells.Item(32, 13).Value = -3875.7256
$ws.Cells.Item(32, 14).Value = -156126.8

$ws.Cells.Item(74, 8).Value = 1102.6389
$ws.Cells.Item(74, 9).Value = 1012.6875
$ws.Cells.Item(74, 10).Value = 1174.6
$ws.Cells.Item(74, 11).Value = 1012.6875
$ws.Cells.Item(74, 12).Value = 1174.6
$ws.Cells.Item(74, 13).Value = -138.6875
$ws.Cells.Item(74, 14).Value = -2922.6

$ws.Cells.Item(77, 8).Value = 1102.6389
$ws.Cells.Item(77, 9).Value = 1012.6875
$ws.Cells.Item(77, 10).Value = 1174.6
$ws.Cells.Item(77, 11).Value = 5063.4375
$ws.Cells.Item(77, 12).Value = 5873
$ws.Cells.Item(77, 13).Value = -695.4375
$ws.Cells.Item(77, 14).Value = -14609

$ws.Cells.Item(132, 8).Value = 3174.7144
$ws.Cells.Item(132, 9).Value = 3217.8276
$ws.Cells.Item(132, 10).Value = 2966.3333
$ws.Cells.Item(132, 11).Value = 9653.4828
$ws.Cells.Item(132, 12).Value = 8898.999899999999
$ws.Cells.Item(132, 13).Value = -7123.4828
$ws.Cells.Item(132, 14).Value = -13958.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 3006.05
$ws.Cells.Item(134, 9).Value = 2766
$ws.Cells.Item(134, 11).Value = 8298
$ws.Cells.Item(134, 13).Value = -5763

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(64, 8).Value = 43208.43
$ws.Cells.Item(64, 10).Value = 43208.43
$ws.Cells.Item(64, 12).Value = 43208.43
$ws.Cells.Item(64, 14).Value = -43704.43

$ws.Cells.Item(67, 8).Value = 43208.43
$ws.Cells.Item(67, 10).Value = 43208.43
$ws.Cells.Item(67, 12).Value = 43208.43
$ws.Cells.Item(67, 14).Value = -44924.43

$ws.Cells.Item(75, 8).Value = 30000
$ws.Cells.Item(75, 10).Value = 30000
$ws.Cells.Item(75, 12).Value = 30000
$ws.Cells.Item(75, 14).Value = -31996

$ws.Cells.Item(78, 8).Value = 30000
$ws.Cells.Item(78, 10).Value = 30000
$ws.Cells.Item(78, 12).Value = 90000
$ws.Cells.Item(78, 14).Value = -99984

$ws.Cells.Item(120, 8).Value = 35266.668
$ws.Cells.Item(120, 10).Value = 35266.668
$ws.Cells.Item(120, 12).Value = 35266.668
$ws.Cells.Item(120, 14).Value = -42524.668

$ws.Cells.Item(121, 8).Value = 46400
$ws.Cells.Item(121, 10).Value = 46400
$ws.Cells.Item(121, 12).Value = 46400
$ws.Cells.Item(121, 14).Value = -49020

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 556506.5
$ws.Cells.Item(107, 9).Value = 438.6111
$ws.Cells.Item(107, 10).Value = 1986395.4
$ws.Cells.Item(107, 11).Value = 1315.8333
$ws.Cells.Item(107, 12).Value = 5959186.199999999
$ws.Cells.Item(107, 13).Value = 604.1667
$ws.Cells.Item(107, 14).Value = -5963026.199999999

$ws.Cells.Item(123, 8).Value = 4253.857
$ws.Cells.Item(123, 9).Value = 2515
$ws.Cells.Item(123, 11).Value = 7545
$ws.Cells.Item(123, 13).Value = -5095

$ws.Cells.Item(131, 8).Value = 848.25
$ws.Cells.Item(131, 9).Value = 364
$ws.Cells.Item(131, 10).Value = 873.7368
$ws.Cells.Item(131, 11).Value = 1092
$ws.Cells.Item(131, 12).Value = 2621.2104
$ws.Cells.Item(131, 13).Value = 3948
$ws.Cells.Item(131, 14).Value = -12701.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 44996.2
$ws.Cells.Item(64, 10).Value = 44996.2
$ws.Cells.Item(64, 12).Value = 44996.2
$ws.Cells.Item(64, 14).Value = -45492.2

$ws.Cells.Item(67, 8).Value = 44996.2
$ws.Cells.Item(67, 10).Value = 44996.2
$ws.Cells.Item(67, 12).Value = 44996.2
$ws.Cells.Item(67, 14).Value = -46712.2

$ws.Cells.Item(75, 8).Value = 18900
$ws.Cells.Item(75, 10).Value = 18900
$ws.Cells.Item(75, 12).Value = 18900
$ws.Cells.Item(75, 14).Value = -20648

$ws.Cells.Item(78, 8).Value = 18900
$ws.Cells.Item(78, 10).Value = 18900
$ws.Cells.Item(78, 12).Value = 56700
$ws.Cells.Item(78, 14).Value = -65436

$ws.Cells.Item(132, 8).Value = 3809.261
$ws.Cells.Item(132, 9).Value = 2464.6365
$ws.Cells.Item(132, 10).Value = 5041.8335
$ws.Cells.Item(132, 11).Value = 7393.9095
$ws.Cells.Item(132, 12).Value = 15125.5005
$ws.Cells.Item(132, 13).Value = -4863.9095
$ws.Cells.Item(132, 14).Value = -20185.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 102250
$ws.Cells.Item(40, 9).Value = 1000000
$ws.Cells.Item(40, 11).Value = 1000000
$ws.Cells.Item(40, 13).Value = -999864

$ws.Cells.Item(55, 8).Value = 943.7308
$ws.Cells.Item(55, 9).Value = 1392.5834
$ws.Cells.Item(55, 10).Value = 559
$ws.Cells.Item(55, 11).Value = 1392.5834
$ws.Cells.Item(55, 12).Value = 559
$ws.Cells.Item(55, 13).Value = -1219.5834
$ws.Cells.Item(55, 14).Value = -905

$ws.Cells.Item(92, 8).Value = 19999.334
$ws.Cells.Item(92, 10).Value = 19999.334
$ws.Cells.Item(92, 12).Value = 19999.334
$ws.Cells.Item(92, 14).Value = -24991.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 10960.5
$ws.Cells.Item(39, 9).Value = 5521
$ws.Cells.Item(39, 11).Value = 5521
$ws.Cells.Item(39, 13).Value = -5108

$ws.Cells.Item(42, 8).Value = 31462.25
$ws.Cells.Item(42, 10).Value = 31462.25
$ws.Cells.Item(42, 12).Value = 31462.25
$ws.Cells.Item(42, 14).Value = -32218.25

$ws.Cells.Item(43, 8).Value = 53000
$ws.Cells.Item(43, 9).Value = 53000
$ws.Cells.Item(43, 11).Value = 53000
$ws.Cells.Item(43, 13).Value = -52851

$ws.Cells.Item(110, 8).Value = 33000
$ws.Cells.Item(110, 10).Value = 33000
$ws.Cells.Item(110, 12).Value = 33000
$ws.Cells.Item(110, 14).Value = -41180

$ws.Cells.Item(116, 8).Value = 49985
$ws.Cells.Item(116, 10).Value = 49985
$ws.Cells.Item(116, 12).Value = 49985
$ws.Cells.Item(116, 14).Value = -59163

Write-Output "Updated 177 cells across 8 sheets (36 rows)."